$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Row 1 (table row 1)
Set-CellText $t 1 1 "69÷8=8, 5"
Set-CellText $t 1 2 "29÷3=9, 2"
Set-CellText $t 1 3 "41÷4=10, 1"
Set-CellText $t 1 4 "60÷3=20, 0"
Set-CellText $t 1 5 "46÷2=23, 0"

# Row 2 (table row 5)
Set-CellText $t 5 1 "35÷2=17, 1"
Set-CellText $t 5 2 "79÷7=11, 2"
Set-CellText $t 5 3 "62÷9=6, 8"
Set-CellText $t 5 4 "74÷4=18, 2"
Set-CellText $t 5 5 "80÷2=40, 0"

# Row 3 (table row 9)
Set-CellText $t 9 1 "93÷2=46, 1"
Set-CellText $t 9 2 "79÷8=9, 7"
Set-CellText $t 9 3 "26÷5=5, 1"
Set-CellText $t 9 4 "67÷5=13, 2"
Set-CellText $t 9 5 "72÷3=24, 0"

# Row 4 (table row 13)
Set-CellText $t 13 1 "70÷3=23, 1"
Set-CellText $t 13 2 "52÷5=10, 2"
Set-CellText $t 13 3 "61÷3=20, 1"
Set-CellText $t 13 4 "88÷8=11, 0"
Set-CellText $t 13 5 "26÷2=13, 0"

# Row 5 (table row 17)
Set-CellText $t 17 1 "48÷4=12, 0"
Set-CellText $t 17 2 "53÷6=8, 5"
Set-CellText $t 17 3 "83÷3=27, 2"
Set-CellText $t 17 4 "62÷7=8, 6"
Set-CellText $t 17 5 "16÷3=5, 1"

Write-Host "Done updating table cells"
